$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill column E (duplicate_image_filename) with "NA" for data rows 2 through 21
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}

# Keep F1 as an empty cell (re-assert to avoid round-trip artifacts)
$ws.Cells.Item(1, 6).Value = ""
